# Updates cryptos list values (Price and Volume(1h) columns) per the
# Aug 19 2023 08:30:22 UTC GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.066.54"
$ws.Range("E2").Value = "  -1.99%  "

$ws.Range("D3").Value = "1.669.24"
$ws.Range("E3").Value = "  -1.50%  "

$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.26%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5115"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.31%  "

$ws.Range("E7").Value = "  -0.22%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2650"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.13%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06419"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.07%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.94"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.12%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07412"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.66%  "

$ws.Range("D12").Value = "1.689.72"
$ws.Range("E12").Value = "  -0.32%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.504"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.20%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5866"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.30%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008547"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.68%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.38"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.78%  "

$ws.Range("D17").Value = "25.976.19"
$ws.Range("E17").Value = "  -2.44%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.956"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.29%  "

$ws.Range("E19").Value = "  -0.11%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.78"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.04%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "190.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.45%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.239"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.55%  "

$ws.Range("E23").Value = "  -0.22%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "145.45"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.54%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.623"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.34%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1204"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.25%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.64"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.02%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06672"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +17.71%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.321"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.16%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.317"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.04%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.539"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.54%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.529"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.98%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.650"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.14%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.020"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.15%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6108"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.96%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.366"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.10%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.712"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.28%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.256"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.70%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01604"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.42%  "

$ws.Range("D40").Value = "1.083.95"
$ws.Range("E40").Value = "  -1.35%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8699"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.47%  "

$ws.Range("E42").Value = "  +0.54%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.71"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.46%  "

$ws.Range("D44").Value = "1.816.92"
$ws.Range("E44").Value = "  -1.98%  "

$ws.Range("E45").Value = "  +5.81%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.49"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.41%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.004"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.25%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.077"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.08%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05222"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.50%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4285"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.03%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.018"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.81%  "
